$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.837750196456909
$ws.Range("B1").Value = 1.826491475105286
$ws.Range("C1").Value = 7.815680027008057
$ws.Range("D1").Value = 0.9835314154624939
$ws.Range("E1").Value = 0.4174738228321075
